$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.149.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.46%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.467.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'558.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.03%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'163.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.00%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '  -1.04%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.467.05"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.05%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '  -4.39%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '  -4.17%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.25%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.920.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.08%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'68.873.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = '  -2.95%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '  -2.52%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.447.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.98%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '  -3.99%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'342.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.39%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '  -2.82%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.78%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'67.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.21%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'3.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.98%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.594.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'8.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.88%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0₃0818"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.21%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.90%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'439.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.33%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '  -4.30%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '  -5.71%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'156.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.19%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'19.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '  -3.85%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '  -3.54%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '  -3.47%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '  -6.11%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '  +1.82%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '  -4.70%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'133.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.35%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '  -4.35%  '
$ws.Range("E51").Style = "Normal"
